$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the "Draudejo grupe" selector value from "IV" to "I"
$ws.Range("C19").Value = "I"

# Update the active cell selection to match the saved view state
$ws.Range("E17").Select()
